$d = $word.ActiveDocument

# --- Part 1: merge the two "FRI Jan 18" / " 11:52:18 IST 2019" runs into one run ---
# A replace-in-place with identical text coalesces the two runs that make up
# that date/time line into a single run, matching the target XML.
$d.Content.Find.Execute("FRI Jan 18 11:52:18 IST 2019", $false, $false, $false, $false, $false, `
    $true, 1, $false, "FRI Jan 18 11:52:18 IST 2019", 2) | Out-Null

# --- Part 2: append a new "purchase details" entry after the last
#     "Amount Received mode ... - CASH AND CLEARD" paragraph ---

# Locate the paragraph that holds the final "CASH AND CLEARD" line (the last
# transaction block in the document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*CASH AND CLEARD*") {
        $target = $p
    }
}

$cursor = $target.Range.Duplicate
$cursor.Collapse(0)

# Blank separator line
$cursor.InsertAfter("`r")
$cursor.Collapse(0)

# "SUN Jan 20" / " 12:51:26 IST 2019"
$cursor.InsertAfter("`rSUN Jan 20")
$cursor.Collapse(0)
$cursor.InsertAfter(" 12:51:26 IST 2019")
$cursor.Collapse(0)

# "Person Name" .... - SRI DHARA
$cursor.InsertAfter("`rPerson Name`t`t`t`t- SRI DHARA")
$cursor.Collapse(0)

# "Bill number" .... - 10298
$cursor.InsertAfter("`rBill number`t`t`t`t- 10298")
$cursor.Collapse(0)

# separator line
$cursor.InsertAfter("`r---------------------------------------------------------------")
$cursor.Collapse(0)

# "Item Name" .... - BEET
$cursor.InsertAfter("`rItem Name`t`t`t`t- BEET")
$cursor.Collapse(0)

# "Number of Pockets" .... - 1
$cursor.InsertAfter("`rNumber of Pockets`t`t`t- 1")
$cursor.Collapse(0)

# "Number of KGs" .... - 57
$cursor.InsertAfter("`rNumber of KGs`t`t`t- 57")
$cursor.Collapse(0)

# "Rate" .... - 20
$cursor.InsertAfter("`rRate`t`t`t`t`t- 20")
$cursor.Collapse(0)

# "Total Price" .... - 1140.0
$cursor.InsertAfter("`rTotal Price`t`t`t`t- 1140.0")
$cursor.Collapse(0)

# "Amount balance" (bold) .... - 1140.0 (bold)
$cursor.InsertAfter("`rAmount balance`t`t`t- 1140.0")
$cursor.Collapse(0)

# Two trailing blank lines
$cursor.InsertAfter("`r`r")
$cursor.Collapse(0)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
